$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tiny floating point precision corrections on existing rows (last-digit
# re-roundings carried over from the source recompute).
$ws.Range("G89").Value = 95935.4559965446
$ws.Range("G90").Value = 77313.25379535164
$ws.Range("G91").Value = 86746.36244998997
$ws.Range("G92").Value = 93772.87719580498

# Append the new monthly data row (2023-12-01).
$ws.Range("A94").Value = 45261
$ws.Range("A94").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B94").Value = 77889.58
$ws.Range("C94").Value = 160452.53
$ws.Range("D94").Value = 240678.8022
$ws.Range("E94").Value = 495798.3177
$ws.Range("F94").Value = 69779.20709593572
$ws.Range("G94").Value = 134645.3736769893
